$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "Liam Livingstone"

function Set-TextValue($range, $value) {
    # Force numeric-looking strings to be stored as text (leading apostrophe),
    # then reset the style so no extra "quote prefix" style sticks to the cell.
    $ws.Range($range).Value = "'" + $value
    $ws.Range($range).Style = "Normal"
}

# Header row
$ws.Range("A1").Value = "matchNo"
$ws.Range("B1").Value = "teamName"
$ws.Range("C1").Value = "batterName"
$ws.Range("D1").Value = "states"
$ws.Range("E1").Value = "runs"
$ws.Range("F1").Value = "balls"
$ws.Range("G1").Value = "fours"
$ws.Range("H1").Value = "sixes"
$ws.Range("I1").Value = "sr"
$ws.Range("J1").Value = "opponentTeamName"
$ws.Range("K1").Value = "venue"
$ws.Range("L1").Value = "date"
$ws.Range("M1").Value = "result"

# Row 2
$ws.Range("A2").Value = "36th"
$ws.Range("B2").Value = "Rajasthan Royals"
$ws.Range("C2").Value = "Liam Livingstone"
$ws.Range("D2").Value = "c †Pant b Avesh Khan"
Set-TextValue "E2" "1"
Set-TextValue "F2" "3"
Set-TextValue "G2" "0"
Set-TextValue "H2" "0"
Set-TextValue "I2" "33.33"
$ws.Range("J2").Value = "Delhi Capitals"
$ws.Range("K2").Value = "Abu Dhabi"
$ws.Range("L2").Value = "September 25"
$ws.Range("M2").Value = "Capitals won by 33 runs"

# Row 3
$ws.Range("A3").Value = "43rd"
$ws.Range("B3").Value = "Rajasthan Royals"
$ws.Range("C3").Value = "Liam Livingstone"
$ws.Range("D3").Value = "c de Villiers b Chahal"
Set-TextValue "E3" "6"
Set-TextValue "F3" "9"
Set-TextValue "G3" "0"
Set-TextValue "H3" "0"
Set-TextValue "I3" "66.66"
$ws.Range("J3").Value = "Royal Challengers Bangalore"
$ws.Range("K3").Value = "Dubai (DSC)"
$ws.Range("L3").Value = "September 29"
$ws.Range("M3").Value = "RCB won by 7 wickets (with 17 balls remaining)"

# Row 4
$ws.Range("A4").Value = "32nd"
$ws.Range("B4").Value = "Rajasthan Royals"
$ws.Range("C4").Value = "Liam Livingstone"
$ws.Range("D4").Value = "c Allen b Arshdeep Singh"
Set-TextValue "E4" "25"
Set-TextValue "F4" "17"
Set-TextValue "G4" "2"
Set-TextValue "H4" "1"
Set-TextValue "I4" "147.05"
$ws.Range("J4").Value = "Punjab Kings"
$ws.Range("K4").Value = "Dubai (DSC)"
$ws.Range("L4").Value = "September 21"
$ws.Range("M4").Value = "Royals won by 2 runs"

# Row 5
$ws.Range("A5").Value = "40th"
$ws.Range("B5").Value = "Rajasthan Royals"
$ws.Range("C5").Value = "Liam Livingstone"
$ws.Range("D5").Value = "c Abdul Samad b Rashid Khan"
Set-TextValue "E5" "4"
Set-TextValue "F5" "6"
Set-TextValue "G5" "0"
Set-TextValue "H5" "0"
Set-TextValue "I5" "66.66"
$ws.Range("J5").Value = "Sunrisers Hyderabad"
$ws.Range("K5").Value = "Dubai (DSC)"
$ws.Range("L5").Value = "September 27"
$ws.Range("M5").Value = "Sunrisers won by 7 wickets (with 9 balls remaining)"

# Row 6
$ws.Range("A6").Value = "54th"
$ws.Range("B6").Value = "Rajasthan Royals"
$ws.Range("C6").Value = "Liam Livingstone"
$ws.Range("D6").Value = "c Tripathi b Ferguson"
Set-TextValue "E6" "6"
Set-TextValue "F6" "6"
Set-TextValue "G6" "1"
Set-TextValue "H6" "0"
Set-TextValue "I6" "100.00"
$ws.Range("J6").Value = "Kolkata Knight Riders"
$ws.Range("K6").Value = "Sharjah"
$ws.Range("L6").Value = "October 07"
$ws.Range("M6").Value = "KKR won by 86 runs"
